$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure account number column (C) keeps numeric-looking values as text
$ws.Range("C2:C5").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "NOUBAIL MOHAMMED"
$ws.Range("B2").Value = "IR801997"
$ws.Range("C2").Value = "007400000313200019604463"
$ws.Range("D2").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E2").Value = "AWB"
$ws.Range("G2").Value = "004/ZZZ/AV2"
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3000

# Row 3
$ws.Range("A3").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("B3").Value = "Q251990"
$ws.Range("C3").Value = "007400000313200019604463"
$ws.Range("D3").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E3").Value = "AWB"
$ws.Range("G3").Value = "004/ZZZ/AV2"
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3000

# Row 4
$ws.Range("A4").Value = "NOUBAIL MOHAMMED"
$ws.Range("B4").Value = "IR801997"
$ws.Range("C4").Value = "007400000313200019604463"
$ws.Range("D4").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E4").Value = "AWB"
$ws.Range("G4").Value = "004/ZZZ/AV2"
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1000

# Row 5
$ws.Range("A5").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("B5").Value = "Q251990"
$ws.Range("C5").Value = "007400000313200019604463"
$ws.Range("D5").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E5").Value = "AWB"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "004/ZZZ/AV2"
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1000

# Row 6 (totals)
$ws.Range("I6").Value = 8000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 8000
